# Update TPM-derived NATMI ligand-receptor metrics (Tgm2-Itgb3) with new TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 30.014089
$ws.Range("H2").Value = 90.042267
$ws.Range("I2").Value = 0.5469606268302545
$ws.Range("J2").Value = 0.5469606268302545
$ws.Range("M2").Value = 0.110552
$ws.Range("N2").Value = 0.331656
$ws.Range("O2").Value = 0.01126249561724847
$ws.Range("P2").Value = 0.01126249561724847
$ws.Range("Q2").Value = 3.318117567128
$ws.Range("R2").Value = 29.863058104152
$ws.Range("S2").Value = 0.006160141662483218
$ws.Range("T2").Value = 0.006160141662483219

$ws.Range("G3").Value = 30.014089
$ws.Range("H3").Value = 90.042267
$ws.Range("I3").Value = 0.5469606268302545
$ws.Range("J3").Value = 0.5469606268302545
$ws.Range("O3").Value = 0.9181055646724333
$ws.Range("P3").Value = 0.9181055646724334
$ws.Range("Q3").Value = 270.4890910636227
$ws.Range("R3").Value = 2434.401819572604
$ws.Range("S3").Value = 0.5021675951495789
$ws.Range("T3").Value = 0.502167595149579

$ws.Range("G4").Value = 30.014089
$ws.Range("H4").Value = 90.042267
$ws.Range("I4").Value = 0.5469606268302545
$ws.Range("J4").Value = 0.5469606268302545
$ws.Range("M4").Value = 0.6933189999999999
$ws.Range("N4").Value = 2.079957
$ws.Range("O4").Value = 0.07063193971031816
$ws.Range("P4").Value = 0.07063193971031817
$ws.Range("Q4").Value = 20.809338171391
$ws.Range("R4").Value = 187.284043542519
$ws.Range("S4").Value = 0.03863289001819236
$ws.Range("T4").Value = 0.03863289001819237

$ws.Range("I5").Value = 0.2046507965132272
$ws.Range("J5").Value = 0.2046507965132272
$ws.Range("M5").Value = 0.110552
$ws.Range("N5").Value = 0.331656
$ws.Range("O5").Value = 0.01126249561724847
$ws.Range("P5").Value = 0.01126249561724847
$ws.Range("Q5").Value = 1.241506919744
$ws.Range("R5").Value = 11.173562277696
$ws.Range("S5").Value = 0.00230487869879663
$ws.Range("T5").Value = 0.00230487869879663

$ws.Range("I6").Value = 0.2046507965132272
$ws.Range("J6").Value = 0.2046507965132272
$ws.Range("O6").Value = 0.9181055646724333
$ws.Range("P6").Value = 0.9181055646724334
$ws.Range("S6").Value = 0.1878910350934397
$ws.Range("T6").Value = 0.1878910350934397

$ws.Range("I7").Value = 0.2046507965132272
$ws.Range("J7").Value = 0.2046507965132272
$ws.Range("M7").Value = 0.6933189999999999
$ws.Range("N7").Value = 2.079957
$ws.Range("O7").Value = 0.07063193971031816
$ws.Range("P7").Value = 0.07063193971031817
$ws.Range("Q7").Value = 7.786022288967999
$ws.Range("R7").Value = 70.07420060071199
$ws.Range("S7").Value = 0.01445488272099085
$ws.Range("T7").Value = 0.01445488272099085

$ws.Range("G8").Value = 13.63015266666667
$ws.Range("H8").Value = 40.890458
$ws.Range("I8").Value = 0.2483885766565184
$ws.Range("J8").Value = 0.2483885766565184
$ws.Range("M8").Value = 0.110552
$ws.Range("N8").Value = 0.331656
$ws.Range("O8").Value = 0.01126249561724847
$ws.Range("P8").Value = 0.01126249561724847
$ws.Range("Q8").Value = 1.506840637605333
$ws.Range("R8").Value = 13.561565738448
$ws.Range("S8").Value = 0.002797475255968625
$ws.Range("T8").Value = 0.002797475255968625

$ws.Range("G9").Value = 13.63015266666667
$ws.Range("H9").Value = 40.890458
$ws.Range("I9").Value = 0.2483885766565184
$ws.Range("J9").Value = 0.2483885766565184
$ws.Range("O9").Value = 0.9181055646724333
$ws.Range("P9").Value = 0.9181055646724334
$ws.Range("Q9").Value = 122.8358990294551
$ws.Range("R9").Value = 1105.523091265096
$ws.Range("S9").Value = 0.2280469344294148
$ws.Range("T9").Value = 0.2280469344294148

$ws.Range("G10").Value = 13.63015266666667
$ws.Range("H10").Value = 40.890458
$ws.Range("I10").Value = 0.2483885766565184
$ws.Range("J10").Value = 0.2483885766565184
$ws.Range("M10").Value = 0.6933189999999999
$ws.Range("N10").Value = 2.079957
$ws.Range("O10").Value = 0.07063193971031816
$ws.Range("P10").Value = 0.07063193971031817
$ws.Range("Q10").Value = 9.450043816700665
$ws.Range("R10").Value = 85.05039435030599
$ws.Range("S10").Value = 0.01754416697113495
$ws.Range("T10").Value = 0.01754416697113495

